# Carga de "requests" al Excel: se sobrescriben las filas 2 y 3 con los
# datos recibidos (carne y celular llegan como texto, y la fecha se
# registra con el formato d/m/yy), igual que hacia el formulario web.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fila 2 (Estudiante 1) ---------------------------------------------
$ws.Range("A2").NumberFormat = "d/m/yy"
$ws.Range("A2").Value = 736858

$ws.Range("B2").Value = "'2011123456"
$ws.Range("B2").Style = "Normal"

$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "'89999999"
$ws.Range("E2").Style = "Normal"

$ws.Range("H2").Value = 1

# --- Fila 3 (Estudiante 2) ---------------------------------------------
$ws.Range("A3").NumberFormat = "d/m/yy"
$ws.Range("A3").Value = 736858

$ws.Range("B3").Value = "'2015567890"
$ws.Range("B3").Style = "Normal"

$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "'59999999"
$ws.Range("E3").Style = "Normal"

$ws.Range("H3").Value = 1
